$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "curso"
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("C4").Value = 55

$ws.Range("F8").Font.Bold = $false

